# Add new timing rows (23-36, 38-41; row 37 intentionally left blank)
# to the "Testing Timings" sheet, growing the used range from A1:I22 to A1:I41.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testing Timings")

$data = @{}
$data[23] = @(3.100000000093132, 0, 0.279, 0.002999999999999947, 0.2570000000000001, 0, 0.05800000000000005, 1.374, 0.2440000000660802)
$data[24] = @(3.400000000023283, 0, 0.262, 0.004000000000000004, 0.256, 0, 0.04700000000000004, 1.6, 0.2239999998892017)
$data[25] = @(2.899999999906868, 0, 0.27, 0.003000000000000003, 0.266, 0, 0.06400000000000006, 1.072, 0.2630000001317967)
$data[26] = @(3.799999999813735, 0, 0.334, 0.004000000000000004, 0.279, 0, 0.05299999999999994, 1.968, 0.2409999999599677)
$data[27] = @(2.900000000139698, 0, 0.262, 0.004000000000000004, 0.261, 0, 0.05800000000000005, 1.087, 0.2180000001317239)
$data[28] = @(2.599999999860302, 0, 0.06, 0.004000000000000004, 0.238, 0, 0.04299999999999998, 1.146, 0.1859999999851425)
$data[29] = @(2.799999999813735, 0, 0.105, 0.004000000000000004, 0.297, 0, 0.04699999999999999, 1.151, 0.2360000000862783)
$data[30] = @(3.199999999953434, 0, 0.292, 0.004000000000000004, 0.273, 0, 0.05700000000000005, 1.093, 0.2830000000813016)
$data[31] = @(3, 0, 0.095, 0.004000000000000004, 0.33, 0, 0.07899999999999996, 1.114, 0.2789999998640269)
$data[32] = @(3.799999999813735, 0, 0.306, 0.003000000000000003, 0.235, 0, 0.06400000000000006, 1.885, 0.2680000000054861)
$data[33] = @(3, 0, 0.304, 0.003000000000000003, 0.274, 0, 0.05100000000000005, 1.196, 0.2379999998538551)
$data[34] = @(3.199999999953434, 0, 0.094, 0.004000000000000004, 0.326, 0, 0.05300000000000005, 1.59, 0.227000000222688)
$data[35] = @(2.699999999953434, 0, 0.097, 0.004000000000000004, 0.325, 0, 0.05099999999999999, 1.083, 0.2269999999953143)
$data[36] = @(3.099999999860302, 0, 0.084, 0.00399999999999999, 0.256, 0, 0.057, 1.461, 0.2239999998892017)
$data[38] = @(3.100000000093132, 0, 0.236, 0.006000000000000005, 0.252, 0, 0.05599999999999994, 1.113, 0.2519999998185085)
$data[39] = @(0.5, 0, 0.236, 0.006000000000000005, 0.252, 0, 0.05599999999999994, 1.113, 44890.77700000007)
$data[40] = @(1.400000000372529, 0, 0.207, 0, 0, 0, 0, 0.04400000000000001, 0.1800000004550384)
$data[41] = @(1.299999999813735, 0, 0.204, 0, 0, 0, 0, 0.03900000000000001, 0.1400000001012813)

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $vals[$c]
    }
}
